# Added author search and sort validation:
#  - refresh the captured JSON total_count test value (91475 -> 91554)
#  - record the current UI selection on the Headers and JSONData sheets
#  - make StatusCode the active tab instead of JSONData

$wb = $excel.ActiveWorkbook

$wsStatus  = $wb.Worksheets.Item("StatusCode")
$wsHeaders = $wb.Worksheets.Item("Headers")
$wsJson    = $wb.Worksheets.Item("JSONData")

# Updated captured value for the author-search / sort validation run.
$wsJson.Range("B2").Value = "91554"

# Headers sheet: selection moves from A1:A4 to B1.
$wsHeaders.Activate() | Out-Null
$wsHeaders.Range("B1").Select() | Out-Null

# JSONData sheet: selection moves to E2 and it is no longer the active tab.
$wsJson.Activate() | Out-Null
$wsJson.Range("E2").Select() | Out-Null

# StatusCode becomes the active sheet/tab.
$wsStatus.Activate() | Out-Null
$wsStatus.Range("A1").Select() | Out-Null
